$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card23")

# Fix header "Event " -> "Event" (drop trailing space)
$ws.Cells.Item(1, 13).Value = "Event"

# Add new header "Correction " (with trailing space) in column N,
# copying the header formatting (bold, border, centered) from M1.
$ws.Cells.Item(1, 14).Value = "Correction "
$ws.Cells.Item(1, 13).Copy()
$ws.Cells.Item(1, 14).PasteSpecial(-4122)

# Fill M2:M12 with "nan" and add an (empty) N column cell for each
# data row, matching the plain formatting already used in column M.
for ($row = 2; $row -le 12; $row++) {
    $ws.Cells.Item($row, 13).Value = "nan"
    $ws.Cells.Item($row, 13).Copy()
    $ws.Cells.Item($row, 14).PasteSpecial(-4122)
}

$excel.CutCopyMode = 0
